# Fix logic problems in the "données13" dataset: columns A and C were
# recomputed for a number of rows (column B, the group sizes, is untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("données13")

$ws.Cells.Item(18, 1).Value = 16.04
$ws.Cells.Item(18, 3).Value = 75

$ws.Cells.Item(20, 1).Value = 17.96
$ws.Cells.Item(20, 3).Value = 101

$ws.Cells.Item(21, 1).Value = 23.45
$ws.Cells.Item(21, 3).Value = 118

$ws.Cells.Item(23, 1).Value = 22.74
$ws.Cells.Item(23, 3).Value = 123

$ws.Cells.Item(24, 1).Value = 42.3
$ws.Cells.Item(24, 3).Value = 125

$ws.Cells.Item(30, 1).Value = 16.96
$ws.Cells.Item(30, 3).Value = 111

$ws.Cells.Item(38, 1).Value = 5.36
$ws.Cells.Item(38, 3).Value = 126

$ws.Cells.Item(39, 1).Value = 18.57
$ws.Cells.Item(39, 3).Value = 95

$ws.Cells.Item(41, 1).Value = 18.47
$ws.Cells.Item(41, 3).Value = 115

$ws.Cells.Item(42, 1).Value = 6.08
$ws.Cells.Item(42, 3).Value = 98

$ws.Cells.Item(43, 1).Value = 7.8100000000000005
$ws.Cells.Item(43, 3).Value = 107

$ws.Cells.Item(44, 1).Value = 15.939999999999998
$ws.Cells.Item(44, 3).Value = 123

$ws.Cells.Item(45, 1).Value = 9.4499999999999993
$ws.Cells.Item(45, 3).Value = 100

$ws.Cells.Item(47, 1).Value = 74.67
$ws.Cells.Item(47, 3).Value = 127

$ws.Cells.Item(48, 1).Value = 49.32
$ws.Cells.Item(48, 3).Value = 127

$ws.Cells.Item(50, 1).Value = 54.510000000000005
$ws.Cells.Item(50, 3).Value = 118

$ws.Cells.Item(53, 1).Value = 12.09
$ws.Cells.Item(53, 3).Value = 118

$ws.Cells.Item(54, 1).Value = 27.91
$ws.Cells.Item(54, 3).Value = 112

$ws.Cells.Item(57, 1).Value = 18.86
$ws.Cells.Item(57, 3).Value = 78

$ws.Cells.Item(59, 1).Value = 49.94
$ws.Cells.Item(59, 3).Value = 116

$ws.Cells.Item(60, 1).Value = 10.08
$ws.Cells.Item(60, 3).Value = 119

$ws.Cells.Item(61, 1).Value = 48.05
$ws.Cells.Item(61, 3).Value = 120

$ws.Cells.Item(62, 1).Value = 9.82
$ws.Cells.Item(62, 3).Value = 69

$ws.Cells.Item(63, 1).Value = 42.01
$ws.Cells.Item(63, 3).Value = 125

$ws.Cells.Item(64, 1).Value = 24.310000000000002
$ws.Cells.Item(64, 3).Value = 127

$ws.Cells.Item(65, 1).Value = 25.28
$ws.Cells.Item(65, 3).Value = 125

$ws.Cells.Item(66, 1).Value = 8.0299999999999994
$ws.Cells.Item(66, 3).Value = 120

$ws.Cells.Item(69, 1).Value = 12.2
$ws.Cells.Item(69, 3).Value = 113

$ws.Cells.Item(70, 1).Value = 60.69
$ws.Cells.Item(70, 3).Value = 110

$ws.Cells.Item(71, 1).Value = 32.269999999999996
$ws.Cells.Item(71, 3).Value = 121

$ws.Cells.Item(74, 1).Value = 7.59
$ws.Cells.Item(74, 3).Value = 103

$ws.Cells.Item(76, 1).Value = 11.76
$ws.Cells.Item(76, 3).Value = 102

$ws.Cells.Item(77, 1).Value = 21.34
$ws.Cells.Item(77, 3).Value = 115

$ws.Cells.Item(79, 1).Value = 9.42
$ws.Cells.Item(79, 3).Value = 113

$ws.Cells.Item(82, 1).Value = 4.08
$ws.Cells.Item(82, 3).Value = 127

$ws.Cells.Item(83, 1).Value = 5.83
$ws.Cells.Item(83, 3).Value = 127

$ws.Cells.Item(85, 1).Value = 15.879999999999999
$ws.Cells.Item(85, 3).Value = 102

$ws.Cells.Item(86, 1).Value = 43.47
$ws.Cells.Item(86, 3).Value = 117

$ws.Cells.Item(87, 1).Value = 51.6
$ws.Cells.Item(87, 3).Value = 99

$ws.Cells.Item(88, 1).Value = 73.429999999999993
$ws.Cells.Item(88, 3).Value = 115

$ws.Cells.Item(90, 1).Value = 15.09
$ws.Cells.Item(90, 3).Value = 108
